$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '26.535.44'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '1.729.55'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -0.75%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '246.13'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.32%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.02%  '
$ws.Cells.Item(7, 5).Value = '  +0.41%  '
$ws.Cells.Item(8, 5).Value = '  -0.91%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.06247'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.11%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '1.727.92'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.84%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.07082'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.55%  '
$ws.Cells.Item(12, 5).Value = '  -1.30%  '
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '4.558'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +1.24%  '
$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '0.6084'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.94%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '77.29'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -0.27%  '
$ws.Cells.Item(16, 5).Value = '  +0.01%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '26.532.29'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -0.19%  '
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '0.000007323'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +6.24%  '
$ws.Cells.Item(19, 2).Value = 'BinanceUSD'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.07%  '
$ws.Cells.Item(20, 5).Value = '  -1.62%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '1.954.24'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.59%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '4.502'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -2.88%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '8.771'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.58%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '5.239'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -1.89%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '137.17'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.99%  '
$ws.Cells.Item(26, 5).Value = '  -0.17%  '
$ws.Cells.Item(27, 5).Value = '  -2.24%  '
$ws.Cells.Item(28, 2).Value = 'BitcoinCash'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '108.46'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.01%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '1.405'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -2.19%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '3.960'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.29%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.08011'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +1.71%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '3.694'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -1.50%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '0.04574'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.51%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -0.04%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '2.620'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.01%  '
$ws.Cells.Item(36, 5).Value = '  +0.44%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.6319'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -1.84%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '0.8910'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -5.84%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '2.003'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.45%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '2.391'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -1.55%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.02%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.01502'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -0.41%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '101.58'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -10.08%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '5.477'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -4.62%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.3890'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.60%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '7.037'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +5.36%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '0.1182'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -2.16%  '
$ws.Cells.Item(48, 5).Value = '  +1.23%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '7.893'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.50%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '30.59'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.49%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '1.251'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -1.74%  '
